$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.1072470739083369;   C = 0.583924513487991;  D = 0.5944036395536997; E = 0.7709757710549013; F = 0.7806389925716273; G = 23 }
    3  = @{ B = 0.573629067650584;    C = 0.8880024064378264; D = 3.979434659499331;  E = 1.994852039500507;  F = 1.953538050233822;  G = 23 }
    4  = @{ B = 0.228737977167174;    C = 1.392066633737142;  D = 7.304711183963507;  E = 2.702722920308981;  F = 2.753551201305114;  G = 23 }
    5  = @{ B = 0.1376631994370348;   C = 1.217887797378426;  D = 7.521227355969471;  E = 2.742485616365101;  F = 2.80058721592678;   G = 23 }
    6  = @{ B = 0.256355043509169;    C = 1.359960095653506;  D = 7.005699073826143;  E = 2.646828115655821;  F = 2.693591441706439;  G = 23 }
    7  = @{ B = 0.09029628155329977;  C = 1.369716599985631;  D = 7.335564264182292;  E = 2.708424683128976;  F = 2.767756307023949;  G = 23 }
    8  = @{ B = 0.1622545047491004;   C = 1.487177404830128;  D = 8.202263652189526;  E = 2.863959436198342;  F = 2.92362286105626;   G = 23 }
    9  = @{ B = 0.08448834823307437;  C = 1.390260216460968;  D = 7.363914233767023;  E = 2.713653300214864;  F = 2.773296749308377;  G = 23 }
    10 = @{ B = 0.1037379453787874;   C = 1.518653016633174;  D = 7.970283297256951;  E = 2.823169016771216;  F = 2.88466950940459;   G = 23 }
    11 = @{ B = -0.05287204651734589; C = 1.303960028208192;  D = 7.512551892089196;  E = 2.740903480987464;  F = 2.801983048518609;  G = 23 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
